$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.408.23"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.32"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "259.92"
$ws.Range("E5").Value = "  -6.79%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5195"
$ws.Range("E7").Value = "  -2.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3253"
$ws.Range("E8").Value = "  -6.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06770"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.53"
$ws.Range("E10").Value = "  -8.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7634"
$ws.Range("E11").Value = "  -5.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07702"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.34"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.02"
$ws.Range("E14").Value = "  -3.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.017"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.86"
$ws.Range("E17").Value = "  -5.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007937"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.428.80"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.071.27"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.555"
$ws.Range("E22").Value = "  -4.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.447"
$ws.Range("E23").Value = "  -6.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.943"
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.40"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.207"
$ws.Range("E26").Value = "  -7.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.643"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.95"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.03"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.159"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.115"
$ws.Range("E31").Value = "  -5.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08702"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04764"
$ws.Range("E33").Value = "  -3.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.121"
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.842"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7010"
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.052"
$ws.Range("E37").Value = "  -7.45%  "
$ws.Range("E38").Value = "  -5.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.169"
$ws.Range("E39").Value = "  -9.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4802"
$ws.Range("E40").Value = "  -7.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "110.91"
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8900"
$ws.Range("E42").Value = "  -7.39%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.627"
$ws.Range("E45").Value = "  -6.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05865"
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4100"
$ws.Range("E47").Value = "  -8.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.936"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.99"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("E50").Value = "  -9.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8839"
$ws.Range("E51").Value = "  -0.45%  "
